$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.191.49'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '2.529.52'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '2.527.01'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.345'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '2.985.47'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '68.025.96'
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '2.510.52'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.62%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('E27').Value = '  +2.35%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('D30').Value = '0.0₃0987'
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '533.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.33'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.52%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.09%  '
$ws.Range('E35').Value = '  -2.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.61'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.69'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.352'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '148.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.557'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0279'
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.72'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('E51').Value = '  -1.19%  '
